$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells (I1, J1). Copy the formatting from the
# existing header cell (H1: bold, centered, bordered) onto the new
# header cells so they reuse the same style as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I (I0) and J (IF) for rows 2-8.
$values = @(
    @(1, 5),
    @(9, 10),
    @(5, 7),
    @(8, 8),
    @(3, 5),
    @(1, 2),
    @(1, 2)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
